# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Brócoli" (Macroferia Regional de Talca)
# at sheet row 151, pushing the existing rows 151-193 down to 152-194.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(151).Insert()

$ws.Cells.Item(151, 1).Value = 5
$ws.Cells.Item(151, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(151, 3).Value = "Maule"
$ws.Cells.Item(151, 4).Value = 44463
$ws.Cells.Item(151, 5).Value = 7
$ws.Cells.Item(151, 6).Value = 100112023
$ws.Cells.Item(151, 7).Value = "Brócoli"
$ws.Cells.Item(151, 8).Value = "Sin especificar"
$ws.Cells.Item(151, 9).Value = "Primera"
$ws.Cells.Item(151, 10).Value = 3000
$ws.Cells.Item(151, 11).Value = 600
$ws.Cells.Item(151, 12).Value = 600
$ws.Cells.Item(151, 13).Value = 600
$ws.Cells.Item(151, 14).Value = "$/unidad"
$ws.Cells.Item(151, 15).Value = "Región del Maule"
$ws.Cells.Item(151, 16).Value = 600
$ws.Cells.Item(151, 17).Value = 1
$ws.Cells.Item(151, 18).Value = "Hortaliza"
